$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 19155.604
$ws.Range("I137").Value = 2135.4546
$ws.Range("J137").Value = 41622.2
$ws.Range("K137").Value = 6406.3638
$ws.Range("L137").Value = 124866.6
$ws.Range("M137").Value = -3856.3638
$ws.Range("N137").Value = -129966.6

$ws.Range("H138").Value = 2565.9246
$ws.Range("J138").Value = 2890.7837
$ws.Range("L138").Value = 8672.3511
$ws.Range("N138").Value = -18952.3511

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16486.844
$ws.Range("I32").Value = 17047.154
$ws.Range("K32").Value = 17047.154
$ws.Range("M32").Value = -16760.154

$ws.Range("H97").Value = 2215.7144
$ws.Range("I97").Value = 2742
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 2742
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -2246
$ws.Range("N97").Value = -1892

$ws.Range("H102").Value = 5349.9
$ws.Range("I102").Value = 4071.2856
$ws.Range("K102").Value = 4071.2856
$ws.Range("M102").Value = -2449.2856

$ws.Range("H110").Value = 1920.5333
$ws.Range("I110").Value = 1618.909
$ws.Range("J110").Value = 2750
$ws.Range("K110").Value = 1618.909
$ws.Range("L110").Value = 2750
$ws.Range("M110").Value = 426.0909999999999
$ws.Range("N110").Value = -6840

$ws.Range("H132").Value = 14012.05
$ws.Range("I132").Value = 1263.9131
$ws.Range("J132").Value = 31259.53
$ws.Range("K132").Value = 3791.7393
$ws.Range("L132").Value = 93778.59
$ws.Range("M132").Value = -1261.7393
$ws.Range("N132").Value = -98838.59

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2071.2856
$ws.Range("I99").Value = 2166.5
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 2166.5
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -668.5
$ws.Range("N99").Value = -4496

$ws.Range("H134").Value = 24851.682
$ws.Range("I134").Value = 26572.537
$ws.Range("K134").Value = 79717.611
$ws.Range("M134").Value = -77182.611

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 81852856
$ws.Range("J6").Value = 255000000
$ws.Range("L6").Value = 255000000
$ws.Range("N6").Value = -255000226

$ws.Range("H16").Value = 1151.5333
$ws.Range("I16").Value = 954.3333
$ws.Range("K16").Value = 954.3333
$ws.Range("M16").Value = -667.3333

$ws.Range("H31").Value = 11777.881
$ws.Range("I31").Value = 22071.79
$ws.Range("J31").Value = 3274.2173
$ws.Range("K31").Value = 22071.79
$ws.Range("L31").Value = 3274.2173
$ws.Range("M31").Value = -21776.79
$ws.Range("N31").Value = -3864.2173

$ws.Range("H34").Value = 11777.881
$ws.Range("I34").Value = 22071.79
$ws.Range("J34").Value = 3274.2173
$ws.Range("K34").Value = 22071.79
$ws.Range("L34").Value = 3274.2173
$ws.Range("M34").Value = -21869.79
$ws.Range("N34").Value = -3678.2173

$ws.Range("H58").Value = 19033.535
$ws.Range("J58").Value = 72509.71000000001
$ws.Range("L58").Value = 72509.71000000001
$ws.Range("N58").Value = -72915.71000000001

$ws.Range("H113").Value = 1151.5333
$ws.Range("I113").Value = 954.3333
$ws.Range("K113").Value = 954.3333
$ws.Range("M113").Value = 1215.6667

$ws.Range("H132").Value = 18684.781
$ws.Range("I132").Value = 25870.904
$ws.Range("J132").Value = 4965.8184
$ws.Range("K132").Value = 77612.712
$ws.Range("L132").Value = 14897.4552
$ws.Range("M132").Value = -75082.712
$ws.Range("N132").Value = -19957.4552

$ws.Range("H134").Value = 651.7560999999999
$ws.Range("I134").Value = 547.5789
$ws.Range("J134").Value = 1971.3334
$ws.Range("K134").Value = 1642.7367
$ws.Range("L134").Value = 5914.0002
$ws.Range("M134").Value = 892.2633000000001
$ws.Range("N134").Value = -10984.0002

$ws.Range("H136").Value = 19033.535
$ws.Range("J136").Value = 72509.71000000001
$ws.Range("L136").Value = 217529.13
$ws.Range("N136").Value = -222629.13

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5946.857
$ws.Range("J62").Value = 6838
$ws.Range("L62").Value = 20514
$ws.Range("N62").Value = -21886

$ws.Range("H65").Value = 5946.857
$ws.Range("J65").Value = 6838
$ws.Range("L65").Value = 61542
$ws.Range("N65").Value = -68406

$ws.Range("H68").Value = 5090.08
$ws.Range("I68").Value = 724.6667
$ws.Range("J68").Value = 6468.6313
$ws.Range("K68").Value = 2174.0001
$ws.Range("L68").Value = 19405.8939
$ws.Range("M68").Value = -1363.0001
$ws.Range("N68").Value = -21027.8939

$ws.Range("H71").Value = 5090.08
$ws.Range("I71").Value = 724.6667
$ws.Range("J71").Value = 6468.6313
$ws.Range("K71").Value = 6522.0003
$ws.Range("L71").Value = 58217.6817
$ws.Range("M71").Value = -2466.0003
$ws.Range("N71").Value = -66329.6817

$ws.Range("H80").Value = 2881.818
$ws.Range("J80").Value = 2881.818
$ws.Range("L80").Value = 8645.454000000002
$ws.Range("N80").Value = -10517.454

$ws.Range("H83").Value = 2881.818
$ws.Range("J83").Value = 2881.818
$ws.Range("L83").Value = 25936.362
$ws.Range("N83").Value = -35296.362

$ws.Range("H131").Value = 801.08
$ws.Range("I131").Value = 420
$ws.Range("J131").Value = 808.8570999999999
$ws.Range("K131").Value = 1260
$ws.Range("L131").Value = 2426.5713
$ws.Range("M131").Value = 3780
$ws.Range("N131").Value = -12506.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1738.2273
$ws.Range("I102").Value = 1662.1666
$ws.Range("K102").Value = 1662.1666
$ws.Range("M102").Value = -40.16660000000002

$ws.Range("H132").Value = 42239.26
$ws.Range("I132").Value = 41094.46
$ws.Range("J132").Value = 44719.668
$ws.Range("K132").Value = 123283.38
$ws.Range("L132").Value = 134159.004
$ws.Range("M132").Value = -120753.38
$ws.Range("N132").Value = -139219.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2914.4285
$ws.Range("I16").Value = 2879.8
$ws.Range("J16").Value = 3001
$ws.Range("K16").Value = 2879.8
$ws.Range("L16").Value = 3001
$ws.Range("M16").Value = -2709.8
$ws.Range("N16").Value = -3341

$ws.Range("H40").Value = 82782.71000000001
$ws.Range("I40").Value = 103452.55
$ws.Range("K40").Value = 103452.55
$ws.Range("M40").Value = -103316.55

$ws.Range("H132").Value = 1591.9756
$ws.Range("I132").Value = 1083.7576
$ws.Range("J132").Value = 3688.375
$ws.Range("K132").Value = 3251.2728
$ws.Range("L132").Value = 11065.125
$ws.Range("M132").Value = -721.2727999999997
$ws.Range("N132").Value = -16125.125

$ws.Range("H136").Value = 27469.7
$ws.Range("I136").Value = 40045
$ws.Range("J136").Value = 4115.5713
$ws.Range("K136").Value = 120135
$ws.Range("L136").Value = 12346.7139
$ws.Range("M136").Value = -117585
$ws.Range("N136").Value = -17446.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1795.625
$ws.Range("I81").Value = 1300.2858
$ws.Range("K81").Value = 2600.5716
$ws.Range("M81").Value = -1539.5716

$ws.Range("H84").Value = 1795.625
$ws.Range("I84").Value = 1300.2858
$ws.Range("K84").Value = 13002.858
$ws.Range("M84").Value = -7698.858

$ws.Range("H132").Value = 1987.0638
$ws.Range("I132").Value = 1751.9706
$ws.Range("J132").Value = 2601.923
$ws.Range("K132").Value = 5255.9118
$ws.Range("L132").Value = 7805.768999999999
$ws.Range("M132").Value = -2725.9118
$ws.Range("N132").Value = -12865.769

$ws.Range("H136").Value = 1155.4193
$ws.Range("J136").Value = 2288.6667
$ws.Range("L136").Value = 6866.000100000001
$ws.Range("N136").Value = -11966.0001
